$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Id values between rows 21 and 22
$ws.Range("A21").Value = 112044197
$ws.Range("A22").Value = 112044200

# Round the Ost (Q) and Nord (R) coordinate values to whole numbers,
# applying the row 21/22 swap that happened alongside the Id swap.
$ws.Range("Q20").Value = 555034
$ws.Range("R20").Value = 6698210

$ws.Range("Q21").Value = 555034
$ws.Range("R21").Value = 6698209

$ws.Range("Q22").Value = 555046
$ws.Range("R22").Value = 6698231

# Remove the Starttid (Z) and Sluttid (AB) values for rows 20-22
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()

$ws.Range("Z21").ClearContents()
$ws.Range("AB21").ClearContents()

$ws.Range("Z22").ClearContents()
$ws.Range("AB22").ClearContents()
